# Update Sprint Backlog to Add Users
# Adjusted Sprint backlog to include officer capabilities to add new users,
# assigned to Ezra. A couple of existing tasks were re-estimated/re-assigned
# to Ryan to make room in the sprint.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing task (row 11): reduce remaining estimate and reassign to Ryan.
$ws.Range("C11").Value = 0.5
$ws.Range("H11").Value = "Ryan"

# Existing task (row 15): reduce remaining estimate and reassign to Ryan.
$ws.Range("C15").Value = 0.5
$ws.Range("H15").Value = "Ryan"

# New backlog row (row 17): officer capability to add new users, assigned to Ezra.
$ws.Range("A17").Value = "Manage User Accounts And Roles"
$ws.Range("B17").Value = "Implement Officer's Adding New Users"
$ws.Range("C17").Value = 2
$ws.Range("H17").Value = "Ezra"

# Widen the "Related User Story" column to fit the new, longer entry.
$ws.Columns.Item(1).ColumnWidth = 29.666666666666668

# Leave the selection on the newly-edited area, matching the author's last cursor position.
$ws.Range("I17").Select() | Out-Null
